$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT without altering its style/number format,
# so that numeric-looking strings (e.g. "318.86") are not coerced into numbers.
function Set-TextValue($range, $value) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

Set-TextValue $ws.Range("D2") "41.774.89"
$ws.Range("E2").Value = "  +0.61%  "
Set-TextValue $ws.Range("D3") "2.475.39"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.13%  "
Set-TextValue $ws.Range("D5") "318.86"
$ws.Range("E5").Value = "  +1.61%  "
Set-TextValue $ws.Range("D6") "93.15"
$ws.Range("E6").Value = "  +2.55%  "
Set-TextValue $ws.Range("D7") "0.553"
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("E8").Value = "  +0.06%  "
Set-TextValue $ws.Range("D9") "0.518"
$ws.Range("E9").Value = "  +1.19%  "
Set-TextValue $ws.Range("D10") "0.0879"
$ws.Range("E10").Value = "  +11.07%  "
Set-TextValue $ws.Range("D11") "33.14"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("E12").Value = "  +0.51%  "
Set-TextValue $ws.Range("D13") "2.858.99"
$ws.Range("E13").Value = "  +0.56%  "
Set-TextValue $ws.Range("D14") "6.93"
$ws.Range("E14").Value = "  +1.20%  "
Set-TextValue $ws.Range("D15") "15.64"
$ws.Range("E15").Value = "  -1.38%  "
Set-TextValue $ws.Range("D16") "2.481.30"
$ws.Range("E16").Value = "  -1.34%  "
Set-TextValue $ws.Range("D17") "0.804"
$ws.Range("E17").Value = "  +3.75%  "
Set-TextValue $ws.Range("D18") "41.716.49"
$ws.Range("E18").Value = "  +0.48%  "
Set-TextValue $ws.Range("D19") "6.50"
$ws.Range("E19").Value = "  +0.16%  "
Set-TextValue $ws.Range("D20") "0.0₃0951"
$ws.Range("E20").Value = "  +1.27%  "
Set-TextValue $ws.Range("D21") "71.15"
$ws.Range("E21").Value = "  +0.08%  "
Set-TextValue $ws.Range("D22") "11.37"
$ws.Range("E22").Value = "  +2.28%  "
Set-TextValue $ws.Range("D23") "242.13"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("E26").Value = "  +0.03%  "
Set-TextValue $ws.Range("D27") "25.22"
$ws.Range("E27").Value = "  +2.89%  "
Set-TextValue $ws.Range("D28") "2.26"
$ws.Range("E28").Value = "  +0.52%  "
Set-TextValue $ws.Range("D29") "9.76"
$ws.Range("E29").Value = "  +1.23%  "
Set-TextValue $ws.Range("D30") "36.97"
$ws.Range("E30").Value = "  +5.00%  "
Set-TextValue $ws.Range("D31") "158.66"
$ws.Range("E31").Value = "  +1.28%  "
Set-TextValue $ws.Range("D32") "5.52"
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("E33").Value = "  +0.00%  "
Set-TextValue $ws.Range("D34") "0.0765"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  -0.56%  "
Set-TextValue $ws.Range("D36") "17.39"
$ws.Range("E36").Value = "  +1.08%  "
Set-TextValue $ws.Range("D37") "1.87"
$ws.Range("E37").Value = "  +5.26%  "
Set-TextValue $ws.Range("D38") "2.92"
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("E40").Value = "  +2.29%  "
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("E42").Value = "  +7.82%  "
Set-TextValue $ws.Range("D43") "1.999.18"
$ws.Range("E43").Value = "  +2.92%  "
Set-TextValue $ws.Range("D44") "19.13"
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("E45").Value = "  +0.96%  "
Set-TextValue $ws.Range("D46") "2.98"
$ws.Range("E46").Value = "  +3.11%  "
Set-TextValue $ws.Range("D47") "9.47"
$ws.Range("E47").Value = "  +4.95%  "
Set-TextValue $ws.Range("D48") "2.716.56"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D49") "98.30"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D50") "76.62"
$ws.Range("E50").Value = "  +7.36%  "
Set-TextValue $ws.Range("D51") "67.36"
$ws.Range("E51").Value = "  +0.48%  "
